$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# --- sheet "year" (sheet3): add new header columns + 3 new data rows ---
# Order below matters: it reproduces the exact shared-string insertion order
# seen in the target workbook (new strings appended starting at index 139).

$ws3.Range("A2").Value = "rou"
$ws3.Range("B2").Value = "all"
$ws3.Range("E1").Value = "select_scaling_year"
$ws3.Range("E2").Value = "1980, 1991,2000,2010"
$ws3.Range("F1").Value = "start_scaling_year"
$ws3.Range("G1").Value = "end_scaling_year"
$ws3.Range("A3").Value = "mkd"
$ws3.Range("E3").Value = "1980, 1992"
$ws3.Range("A4").Value = "idn"
$ws3.Range("H4").Value = "Eliminate jump in CO emissions present in EDGAR data"
$ws3.Range("H1").Value = "Comment"
$ws3.Range("H2").Value = "Reduce jumps in emissions"

# Remaining cells reuse existing shared strings (NA, all, 1A2) or are numeric.
$ws3.Range("C2").Value = "NA"
$ws3.Range("D2").Value = "NA"
$ws3.Range("F2").Value = "NA"
$ws3.Range("G2").Value = "NA"

$ws3.Range("B3").Value = "all"
$ws3.Range("C3").Value = "NA"
$ws3.Range("D3").Value = "NA"
$ws3.Range("F3").Value = "NA"
$ws3.Range("G3").Value = "NA"
$ws3.Range("H3").Value = "Reduce jumps in emissions"

$ws3.Range("B4").Value = "1A2"
$ws3.Range("C4").Value = "NA"
$ws3.Range("D4").Value = "NA"
$ws3.Range("E4").Value2 = 1990
$ws3.Range("F4").Value2 = 2010
$ws3.Range("G4").Value = "NA"

# Column E is widened in the edited workbook.
$ws3.Columns.Item(5).ColumnWidth = 17.666666666666668

# --- view-state updates ---
# The "map" sheet (sheet1) had its bottom-right pane scrolled/selected before
# the author moved on to the "year" sheet (sheet3), which becomes the tab
# that is active/selected when the workbook is saved.
$ws1.Activate()
$ws1.Range("D42").Select()

$ws3.Activate()
$ws3.Range("H3").Select()
